$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New error code rows (10024-10029) that only carry a code, no message/level yet
$plainCodes = @(10024, 10025, 10026, 10027, 10028, 10029)
$r = 26
foreach ($code in $plainCodes) {
    $ws.Cells.Item($r, 1).Value = $code
    $r++
}

# New user-management related error codes (10030-10032) with message + level
$ws.Range("A32").Value = 10030
$ws.Range("B32").Value = "message_10030_user_record_created_successfully"
$ws.Range("D32").Value = "Success"

$ws.Range("A33").Value = 10031
$ws.Range("B33").Value = "message_10031_user_record_updated_successfully"
$ws.Range("D33").Value = "Success"

$ws.Range("A34").Value = 10032
$ws.Range("B34").Value = "message_10032_user_record_deleted_successfully"
$ws.Range("D34").Value = "Success"

# Trailing plain code rows (10033-10035)
$ws.Range("A35").Value = 10033
$ws.Range("A36").Value = 10034
$ws.Range("A37").Value = 10035

# Match the saved view/selection state: scrolled so row 7 is at top, cell B35 selected
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("B35").Select()
